# Update countries & provincias Spain
# Refresh COVID-19 country stats on sheet "Pais" and bump the "datos
# actualizados" timestamp. A few countries (Paises Bajos / Croacia) pulled
# ahead of their neighbours in total cases, so those rows swap places too -
# but since the sheet is sorted by "Casos totales" desc, we just need to
# write the new totals into the right rows directly.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# --- timestamp banner -------------------------------------------------
$ws.Range("A1").Value = "Datos actualizados a 28 de Agosto de 2020 a las 15:00"

# --- helper: write a full data row (Casos totales..Muertes) -----------
function Set-Row {
    param(
        [int]$Row,
        [string]$Pais,
        $TotalCases,
        $NuevosCasos,
        $CasosActivos,
        $Recuperados,
        $CasosCriticos,
        $MuertesHoy,
        $Muertes
    )
    $ws.Range("A$Row").Value = $Pais
    $ws.Range("B$Row").Value = $TotalCases
    $ws.Range("C$Row").Value = $NuevosCasos
    $ws.Range("D$Row").Value = $CasosActivos
    $ws.Range("E$Row").Value = $Recuperados
    $ws.Range("F$Row").Value = $CasosCriticos
    $ws.Range("G$Row").Value = $MuertesHoy
    $ws.Range("H$Row").Value = $Muertes
}

# Estados Unidos
Set-Row 4 "Estados Unidos" 6049440 2806 3348784 2515729 0 131 184927

# India
Set-Row 6 "India" 3392367 7792 2585037 745605 0 31 61725

# Arabia Saudita
Set-Row 17 "Arabia Saudita" 312924 1069 287403 21708 0 28 3813

# Suecia
Set-Row 40 "Suecia" 83958 0 0 0 0 5 5821

# Paises Bajos now edges out Emiratos Arabes Unidos (rows swap)
Set-Row 45 "Paises Bajos" 69131 507 0 0 0 2 6220
Set-Row 46 "Emiratos Arabes Unidos" 68901 390 59861 8661 0 1 379

# Uzbekistan
Set-Row 62 "Uzbekistan" 40720 273 37873 2545 0 4 302

# Dinamarca
Set-Row 80 "Dinamarca" 16700 73 14877 1199 0 0 624

# Croacia jumps ahead of Grecia / Malasia / Guinea (rows shift down one)
Set-Row 92 "Croacia" 9549 357 6809 2560 0 3 180
Set-Row 93 "Grecia" 9531 0 3804 5473 0 0 254
Set-Row 94 "Malasia" 9306 10 9030 151 0 0 125
Set-Row 95 "Guinea" 9213 0 8180 975 0 0 58

# Togo
Set-Row 156 "Togo" 1341 15 968 346 0 0 27

# Islas Feroe (only Casos activos / Recuperados changed)
$ws.Range("D179").Value = 361
$ws.Range("E179").Value = 50
